# Automatische test-sync: 2025-08-13 21:56:50
# Appends a new "Demo inplannen" log row to the Logs sheet, extends the
# conditional-formatting ranges that covered the data rows, and bumps the
# matching tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Find the next free row right after the current data (row 20 -> 21).
$newRow = $ws.UsedRange.Rows.Count + 1

$ws.Cells.Item($newRow, 1).Value  = "Demo inplannen"
$ws.Cells.Item($newRow, 2).Value  = "klantenservice@testbedrijf123.nl"
$ws.Cells.Item($newRow, 3).Value  = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item($newRow, 4).Value  = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Cells.Item($newRow, 6).Value  = "2025-08-13 21:55:58"
$ws.Cells.Item($newRow, 7).Value  = "Nee"
$ws.Cells.Item($newRow, 8).Value  = "Ja"
$ws.Cells.Item($newRow, 9).Value  = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Extend the five conditional-formatting blocks (D, G, H, I, J) so they
# keep covering the data rows, now including the freshly added row.
$oldLastRow = $newRow - 1

$columns = "D", "G", "H", "I", "J"
foreach ($col in $columns) {
    $oldRange = $ws.Range("$col" + "2:" + "$col" + "$oldLastRow")
    $newRange = $ws.Range("$col" + "2:" + "$col" + "$newRow")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard tally for the "Intern verzoek / Actie voor
# medewerker" category to reflect the newly logged row.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = $dash.Cells.Item(2, 2).Value() + 1
